# BABA.xlsx update — refresh the price input and the model's discount rate,
# then leave the selection/active sheet the way the author left it.

$wb = $excel.ActiveWorkbook

$wsMain  = $wb.Worksheets.Item("Main")
$wsModel = $wb.Worksheets.Item("Model")

# --- Main!D3 (Price): 121 -> 117, shown with 2 decimals now ---
$wsMain.Range("D3").Value = 117
$wsMain.Range("D3").NumberFormat = "#,##0.00"

# --- Model!S21 (Discount rate): 6% -> 4% ---
$wsModel.Range("S21").Value = 0.04

# --- Model!S25 (Price) picks up the same 2-decimal format as Main!D3 ---
$wsModel.Range("S25").NumberFormat = "#,##0.00"

# --- Leave the cursor/selection where the author left it: Model!B36, then
#     switch back to the Main tab with D4 selected ---
$wsModel.Range("B36").Select()
$wsMain.Activate()
$wsMain.Range("D4").Select()
